$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Edit 1 - Slide 13, shape "TextBox 42": "but K" -> "but S"
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$shp13 = $s13.Shapes.Item("TextBox 42")
$tr13 = $shp13.TextFrame.TextRange
$found13 = $tr13.Find("but K")
$found13.Text = "but S"

# ---------------------------------------------------------------------------
# Edit 2 - Slide 20, shape "직사각형 15": remove the " -ansi" fragment so the
# line reads "... -Wall -Werror -pedantic -std=c99 " instead of
# "... -Wall -Werror -ansi -pedantic -std=c99 ".
# ---------------------------------------------------------------------------
$s20 = $p.Slides.Item(20)
$shp20 = $s20.Shapes.Item("직사각형 15")
$tr20 = $shp20.TextFrame.TextRange
$found20 = $tr20.Find(" -ansi")
$found20.Text = ""

# ---------------------------------------------------------------------------
# Edit 3 - Slide 27, shape "직사각형 3": same fix in the shell-script textbox,
# turning "gcc -Wall -Werror -ansi -pedantic -std=c99 "$@"" into
# "gcc -Wall -Werror -pedantic -std=c99 "$@"".
# ---------------------------------------------------------------------------
$s27 = $p.Slides.Item(27)
$shp27 = $s27.Shapes.Item("직사각형 3")
$tr27 = $shp27.TextFrame.TextRange
$found27 = $tr27.Find("ansi -")
$found27.Text = ""
